# Add the required "Experimental" boolean element (row 7, column B) on the
# Metadata sheet, and bump the Date value (row 8, column B) to reflect the
# regenerated valueset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$dst = $ws.Range("B7")
$src = $ws.Range("B13")

# Write the literal text "true" (quote-prefixed so the engine stores it as
# text instead of inferring a Boolean), then restore the plain un-prefixed
# cell style that every other value cell in this column already uses.
$dst.Value = "'true"
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
